# Shopping slots export: add "대행사" (Agency) and "상태" (Status) columns,
# and trim the sheet back down to just the header row (remove the two
# sample data rows that used to ship with the export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B for "대행사" (Agency). Everything from the old B
# ("스토어 타입") onward shifts one column to the right.
$ws.Columns("B:B").Insert()

# Insert a new column L for "상태" (Status), just before the last column
# ("슬롯 단가", which is now sitting in column L and shifts to M).
$ws.Columns("L:L").Insert()

# Fill in the headers for the two newly inserted columns.
$ws.Range("B1").Value = "대행사"
$ws.Range("L1").Value = "상태"

# Drop the two sample data rows - only the header row remains.
$ws.Rows("2:3").Delete()
